$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the price/volume columns stay formatted as Text so the
# numeric-looking strings are not auto-converted to numbers/percentages.
$ws.Range("D2:E26").NumberFormat = "@"
$ws.Range("D38:E51").NumberFormat = "@"

$ws.Range("D2").Value = "326.62"
$ws.Range("E2").Value = "-1.17%"
$ws.Range("D3").Value = "44.20"
$ws.Range("E3").Value = "0.38%"
$ws.Range("D4").Value = "5.245"
$ws.Range("E4").Value = "-5.21%"
$ws.Range("D5").Value = "0.08315"
$ws.Range("E5").Value = "2.50%"
$ws.Range("D6").Value = "1.930"
$ws.Range("E6").Value = "-6.24%"
$ws.Range("D7").Value = "0.9703"
$ws.Range("E7").Value = "-0.39%"
$ws.Range("D8").Value = "2.527"
$ws.Range("E8").Value = "-3.25%"
$ws.Range("D9").Value = "0.1134"
$ws.Range("E9").Value = "2.95%"
$ws.Range("D10").Value = "0.1888"
$ws.Range("E10").Value = "-0.28%"
$ws.Range("D11").Value = "0.09622"
$ws.Range("E11").Value = "-3.46%"
$ws.Range("D12").Value = "0.04620"
$ws.Range("E12").Value = "-2.24%"
$ws.Range("D13").Value = "0.1057"
$ws.Range("E13").Value = "0.19%"
$ws.Range("D14").Value = "0.001289"
$ws.Range("E14").Value = "1.76%"
$ws.Range("D15").Value = "0.006124"
$ws.Range("E15").Value = "0.94%"
$ws.Range("D16").Value = "3.397"
$ws.Range("E16").Value = "1.65%"
$ws.Range("D17").Value = "4.436"
$ws.Range("E17").Value = "0.00%"
$ws.Range("D18").Value = "0.3346"
$ws.Range("E18").Value = "0.02%"
$ws.Range("D19").Value = "8.755"
$ws.Range("E19").Value = "-13.71%"
$ws.Range("D20").Value = "0.1373"
$ws.Range("E20").Value = "-0.66%"
$ws.Range("D21").Value = "0.2582"
$ws.Range("E21").Value = "0.41%"
$ws.Range("D22").Value = "0.04151"
$ws.Range("E22").Value = "0.94%"
$ws.Range("D23").Value = "0.001234"
$ws.Range("E23").Value = "-5.41%"
$ws.Range("D24").Value = "0.004409"
$ws.Range("E24").Value = "0.83%"
$ws.Range("D25").Value = "0.0001304"
$ws.Range("E25").Value = "1.94%"
$ws.Range("D26").Value = "0.0002989"
$ws.Range("E26").Value = "-20.02%"
$ws.Range("D38").Value = "0.02683"
$ws.Range("E38").Value = "0.25%"
$ws.Range("D39").Value = "0.05545"
$ws.Range("E39").Value = "-1.50%"
$ws.Range("D40").Value = "0.007841"
$ws.Range("E40").Value = "3.21%"
$ws.Range("D41").Value = "0.1407"
$ws.Range("E41").Value = "-0.57%"
$ws.Range("D42").Value = "0.007352"
$ws.Range("E42").Value = "-2.59%"
$ws.Range("D43").Value = "0.002123"
$ws.Range("E43").Value = "8.48%"
$ws.Range("D44").Value = "0.007867"
$ws.Range("E44").Value = "-5.39%"
$ws.Range("D45").Value = "0.3497"
$ws.Range("D46").Value = "0.00006848"
$ws.Range("E46").Value = "-2.65%"
$ws.Range("D47").Value = "0.00000000753"
$ws.Range("E47").Value = "0.40%"
$ws.Range("D48").Value = "0.003496"
$ws.Range("E48").Value = "-1.19%"
$ws.Range("D49").Value = "0.003542"
$ws.Range("E49").Value = "40.65%"
$ws.Range("D50").Value = "0.00002107"
$ws.Range("E50").Value = "0.40%"
$ws.Range("D51").Value = "0.0002007"
$ws.Range("E51").Value = "0.40%"
